$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the last existing row (67) into the new row (68) first so that the
# new row inherits the same cell formatting/styles (e.g. bold/border style
# on column A, and the date number-format on column E).
$ws.Range("A67:V67").Copy($ws.Range("A68:V68"))

# Now overwrite the copied values with the new match's data.
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = "croatia"
$ws.Range("C68").Value = "hnl"
$ws.Range("D68").Value = "2023-2024"
$ws.Range("E68").Value = 45238.70833333334
$ws.Range("F68").Value = "Osijek"
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = "Hajduk Split"
$ws.Range("I68").Value = 1
$ws.Range("J68").Value = 2.34
$ws.Range("K68").Value = "30/07/2023 21:12"
$ws.Range("L68").Value = 3.11
$ws.Range("M68").Value = "08/11/2023 16:59"
$ws.Range("N68").Value = 3.22
$ws.Range("O68").Value = "30/07/2023 21:12"
$ws.Range("P68").Value = 3.23
$ws.Range("Q68").Value = "08/11/2023 16:59"
$ws.Range("R68").Value = 3.01
$ws.Range("S68").Value = "30/07/2023 21:12"
$ws.Range("T68").Value = 2.42
$ws.Range("U68").Value = "08/11/2023 16:56"
$ws.Range("V68").Value = "https://www.betexplorer.com/football/croatia/hnl/osijek-hajduk-split/UyrQj0Pu/"
